# Applies the "Saldo" export update:
#   - adds several new account rows (keeping the Saldo-descending sort order)
#   - bumps BRUNO's (004584982) balance from 5000 to 5490.84
#   - moves SANDRA (004472760) from a 24.96 balance row further down the
#     sheet to a new row near the top with a 1123.13 balance
#   - adds a final MARIANA (004525587) row with a negative balance
#
# Work bottom-to-top so row numbers captured from the original layout stay
# valid for every subsequent operation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: this interpreter only supports positional arguments for user
# functions (named parameters on custom functions are not resolved), so
# call Insert-DataRow positionally everywhere below.
#
# The account numbers ("Conta") are zero-padded numeric-looking strings
# (e.g. "004001621"); writing them as a bare Value lets Excel coerce them
# to a real number and drop the leading zeros. Prefixing with a leading
# apostrophe forces text, exactly like typing '004001621 into the cell.
function Insert-DataRow {
    param($RowIndex, $Conta, $Nome, $Saldo)
    $ws.Rows.Item($RowIndex).Insert()
    $ws.Cells.Item($RowIndex, 1).Value = "'" + $Conta
    $ws.Cells.Item($RowIndex, 2).Value = $Nome
    $ws.Cells.Item($RowIndex, 3).Value = $Saldo
}

# 7) New row for MARIANA, right before the trailing blank row (row 254).
Insert-DataRow 254 "004525587" "MARIANA" -0.08

# 6) Remove the old SANDRA row (balance 24.96) further down the sheet.
$ws.Rows.Item(157).Delete()

# 5) New row for PAULA, just above RODRIGO (row 16).
Insert-DataRow 16 "004503374" "PAULA" 934.89

# 4) New row for SANDRA with her updated balance, just above BLUEMETRIX (row 14).
Insert-DataRow 14 "004472760" "SANDRA" 1123.13

# 3) New row for PATRICIA, just above GUSTAVO (row 8).
Insert-DataRow 8 "004267044" "PATRICIA" 4392.67

# 2) BRUNO's (004584982) balance changes from 5000 to 5490.84 (row 7).
$ws.Cells.Item(7, 3).Value = 5490.84

# 1) New rows for DANIELA and ASSAKO, just above CLAUDIA (row 3).
Insert-DataRow 3 "004001621" "DANIELA" 58543.89
Insert-DataRow 4 "004450724" "ASSAKO" 57104.67
